$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the header style from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$null = $ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new I/J column data for rows 2-19
$data = @(
    @(5, 6),
    @(8, 8),
    @(5, 5),
    @(8, 9),
    @(7, 7),
    @(8, 9),
    @(6, 8),
    @(7, 9),
    @(5, 6),
    @(3, 7),
    @(1, 5),
    @(5, 6),
    @(1, 4),
    @(1, 2),
    @(1, 1),
    @(3, 4),
    @(5, 5),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
